$d = $word.ActiveDocument

# 1. Remove the title text ("React Native Birthday - Push Notifications")
#    from the first paragraph while keeping its paragraph formatting
#    (alignment + the cached run formatting stored on the paragraph mark).
$d.Content.Find.Execute("React Native Birthday - Push Notifications", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Insert a brand-new (empty) paragraph right after it.
$endOfDoc = $d.Range($d.Content.End, $d.Content.End)
$endOfDoc.InsertParagraphAfter()

$newPara = $d.Paragraphs(2)

# 3. Turn on underline for that new paragraph. Directly setting the
#    underline on the still-empty paragraph mark would strand a leftover
#    empty run, so briefly insert a placeholder character, format it, and
#    then remove the character again -- this leaves the underline setting
#    cached on the paragraph mark (w:pPr/w:rPr) with no stray w:r behind.
$placeholder = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$placeholder.InsertAfter("X")

$newPara = $d.Paragraphs(2)
$newPara.Range.Font.Underline = 1

$charRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$charRange.Delete()
